$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 91
$ws.Range("D3").Value = 132
$ws.Range("J3").Value = 225
$ws.Range("K4").Value = 23
$ws.Range("B6").Value = 368
$ws.Range("C6").Value = 471
$ws.Range("D6").Value = 407
$ws.Range("E6").Value = 461
$ws.Range("F6").Value = 515
$ws.Range("G6").Value = 431
$ws.Range("I6").Value = 495
$ws.Range("J6").Value = 411
$ws.Range("K6").Value = 499
$ws.Range("B7").Value = 493
$ws.Range("C7").Value = 625
$ws.Range("D7").Value = 636
$ws.Range("E7").Value = 684
$ws.Range("F7").Value = 747
$ws.Range("G7").Value = 660
$ws.Range("I7").Value = 827
$ws.Range("J7").Value = 778
$ws.Range("K7").Value = 877

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("C6").Value = 34
$ws.Range("I6").Value = 32
$ws.Range("K6").Value = 25
$ws.Range("C7").Value = 39
$ws.Range("I7").Value = 48
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("G6").Value = 15
$ws.Range("G7").Value = 27

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 15

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("D3").Value = 17
$ws.Range("D7").Value = 45

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("D4").Value = 6
$ws.Range("D5").Value = 10

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("D10").Value = 3
$ws.Range("F19").Value = 23
$ws.Range("D28").Value = 45
$ws.Range("C32").Value = 39
$ws.Range("I32").Value = 48
$ws.Range("K32").Value = 45
$ws.Range("G36").Value = 27
$ws.Range("J42").Value = 13
$ws.Range("D47").Value = 14
$ws.Range("J47").Value = 16
$ws.Range("E53").Value = 81
$ws.Range("F53").Value = 79
$ws.Range("J53").Value = 120
$ws.Range("K61").Value = 4
$ws.Range("C63").Value = 7
$ws.Range("K76").Value = 29
$ws.Range("B80").Value = 15
$ws.Range("D82").Value = 10
$ws.Range("C85").Value = 15
$ws.Range("B98").Value = 493
$ws.Range("C98").Value = 625
$ws.Range("D98").Value = 636
$ws.Range("E98").Value = 684
$ws.Range("F98").Value = 747
$ws.Range("G98").Value = 660
$ws.Range("I98").Value = 827
$ws.Range("J98").Value = 778
$ws.Range("K98").Value = 877

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("F2").Value = 7
$ws.Range("E6").Value = 63
$ws.Range("F6").Value = 58
$ws.Range("J6").Value = 60
$ws.Range("E7").Value = 81
$ws.Range("F7").Value = 79
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 29

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("H4").Value = 11
$ws.Range("H5").Value = 13

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 3

$ws = $wb.Worksheets.Item("New City")
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 7

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 23

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("C4").Value = 12
$ws.Range("C5").Value = 15

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 11
$ws.Range("D5").Value = 12
$ws.Range("D6").Value = 14
$ws.Range("J6").Value = 16
